{"js": "// Update the date line at the top of the document.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph\n  .getRange()\n  .insertText(\"2024-08-10 Saturday\", Word.InsertLocation.replace);\n\n// Update the division problems in the table, cell by cell\n// (rowIndex, cellIndex are 0-based).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  [0, 0, \"78\u00f74=\"],\n  [0, 1, \"58\u00f74=\"],\n  [0, 2, \"24\u00f78=\"],\n  [0, 3, \"35\u00f77=\"],\n  [0, 4, \"53\u00f73=\"],\n\n  [4, 0, \"62\u00f74=\"],\n  [4, 1, \"18\u00f73=\"],\n  [4, 2, \"19\u00f78=\"],\n  [4, 3, \"12\u00f77=\"],\n  [4, 4, \"30\u00f73=\"],\n\n  [8, 0, \"63\u00f77=\"],\n  [8, 1, \"34\u00f77=\"],\n  [8, 2, \"10\u00f78=\"],\n  [8, 3, \"27\u00f79=\"],\n  [8, 4, \"95\u00f74=\"],\n\n  [12, 0, \"71\u00f74=\"],\n  [12, 1, \"53\u00f77=\"],\n  [12, 2, \"69\u00f74=\"],\n  [12, 3, \"85\u00f72=\"],\n  [12, 4, \"53\u00f72=\"],\n\n  [16, 0, \"89\u00f73=\"],\n  [16, 1, \"96\u00f73=\"],\n  // [16, 2] (\"68\u00f72=\") is unchanged.\n  [16, 3, \"66\u00f79=\"],\n  [16, 4, \"14\u00f77=\"],\n];\n\nfor (const [rowIndex, cellIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, cellIndex);\n  const paragraph = cell.body.paragraphs.getFirst();\n  paragraph.getRange().insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line at the top of the document.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"2024-08-09 Friday\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"2024-08-10 Saturday\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# Update the division problems in the table, cell by cell (row, col are 1-based).\n$t = $d.Tables(1)\n\n$t.Cell(1,1).Range.Text = \"78\u00f74=\"\n$t.Cell(1,2).Range.Text = \"58\u00f74=\"\n$t.Cell(1,3).Range.Text = \"24\u00f78=\"\n$t.Cell(1,4).Range.Text = \"35\u00f77=\"\n$t.Cell(1,5).Range.Text = \"53\u00f73=\"\n\n$t.Cell(5,1).Range.Text = \"62\u00f74=\"\n$t.Cell(5,2).Range.Text = \"18\u00f73=\"\n$t.Cell(5,3).Range.Text = \"19\u00f78=\"\n$t.Cell(5,4).Range.Text = \"12\u00f77=\"\n$t.Cell(5,5).Range.Text = \"30\u00f73=\"\n\n$t.Cell(9,1).Range.Text = \"63\u00f77=\"\n$t.Cell(9,2).Range.Text = \"34\u00f77=\"\n$t.Cell(9,3).Range.Text = \"10\u00f78=\"\n$t.Cell(9,4).Range.Text = \"27\u00f79=\"\n$t.Cell(9,5).Range.Text = \"95\u00f74=\"\n\n$t.Cell(13,1).Range.Text = \"71\u00f74=\"\n$t.Cell(13,2).Range.Text = \"53\u00f77=\"\n$t.Cell(13,3).Range.Text = \"69\u00f74=\"\n$t.Cell(13,4).Range.Text = \"85\u00f72=\"\n$t.Cell(13,5).Range.Text = \"53\u00f72=\"\n\n$t.Cell(17,1).Range.Text = \"89\u00f73=\"\n$t.Cell(17,2).Range.Text = \"96\u00f73=\"\n$t.Cell(17,4).Range.Text = \"66\u00f79=\"\n$t.Cell(17,5).Range.Text = \"14\u00f77=\"\n"}
